$d = $word.ActiveDocument
$brk = [string][char]11

# Find the "Event not competed in by this team" paragraph that
# immediately follows the "Girl's 3200:" header paragraph - this is
# where the two new event blocks (Girl's HH: and Girl's 300H:) need
# to be inserted.
$idx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "Girl's 3200:") {
        $idx = $i
        break
    }
}
$anchorPara = $d.Paragraphs.Item($idx + 1)

# Insert a new, empty paragraph right after $afterPara and return it.
function New-ParagraphAfter($afterPara) {
    $afterPara.Range.InsertParagraphAfter()
    return $afterPara.Next()
}

# Fill a freshly-created empty paragraph with bold text (used for the
# "Girl's <Event>:" header lines, which consist of two leading line
# breaks, the bold label, and a trailing line break).
function Set-BoldHeaderText($para, $text, $doc) {
    $r = $para.Range
    $r.Collapse(1)
    $start = $r.Start
    $r.InsertAfter($text)
    $tail = $doc.Range($start, $start + $text.Length)
    $tail.Bold = $true
}

# Fill a freshly-created empty paragraph with plain (non-bold) text
# (used for the "Event not competed in by this team" lines).
function Set-PlainText($para, $text, $doc) {
    $r = $para.Range
    $r.Collapse(1)
    $start = $r.Start
    $r.InsertAfter($text)
    $tail = $doc.Range($start, $start + $text.Length)
    $tail.Bold = $false
}

$notCompetedText = "Event not competed in by this team"
$hhHeaderText = $brk + $brk + "Girl's HH:" + $brk
$hurdleHeaderText = $brk + $brk + "Girl's 300H:" + $brk

$hhHeader = New-ParagraphAfter $anchorPara
Set-BoldHeaderText $hhHeader $hhHeaderText $d

$hhEvent = New-ParagraphAfter $hhHeader
Set-PlainText $hhEvent $notCompetedText $d

$hurdleHeader = New-ParagraphAfter $hhEvent
Set-BoldHeaderText $hurdleHeader $hurdleHeaderText $d

$hurdleEvent = New-ParagraphAfter $hurdleHeader
Set-PlainText $hurdleEvent $notCompetedText $d
